$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.599.48"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "2.440.40"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.29"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.59"
$ws.Range("E6").Value = "  +9.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "2.444.51"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  +5.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.124"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "2.880.35"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "56.971.07"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.86"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "2.451.58"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.96"
$ws.Range("E20").Value = "  +4.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.00"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.66"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.403"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").Value = "2.566.16"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").Value = "0.0₃0804"
$ws.Range("E30").Value = "  +6.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.02"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.51"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.05"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.886"
$ws.Range("E36").Value = "  +6.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.14"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.72"
$ws.Range("E38").Value = "  +5.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.39"
$ws.Range("E39").Value = "  +9.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.04"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0557"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.604"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  +7.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.84"
$ws.Range("E46").Value = "  +4.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "263.71"
$ws.Range("E47").Value = "  +3.76%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.26"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.67"
$ws.Range("E50").Value = "  +4.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("E51").Value = "  +27.23%  "
